$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before D. This pushes the old D/F/G/H columns
#    (notes / microSD_id / formatted / per-row notes) one slot to the right,
#    matching the new "set_local_time" column that the author added to the
#    Audiomoth table.
# ---------------------------------------------------------------------------
$ws.Columns("D").Insert()

# New column header for the Audiomoth table.
$ws.Range("D1").Value = "set_local_time"

# ---------------------------------------------------------------------------
# 2. Re-type the SongMini header row (row 28). After the column insert its
#    old contents landed on B/C/E/F; overwrite with the final header layout.
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = "firmware"
$ws.Range("C28").Value = "configured"
$ws.Range("D28").Value = "set_local_time"
$ws.Range("E28").Value = "notes"
$ws.Range("F28").ClearContents()
$ws.Range("G28").Value = "SD"
$ws.Range("H28").Value = "formatted"

# ---------------------------------------------------------------------------
# 3. Clear out the old SongMini data rows (29:38) which only held leftover
#    "exFAT" filler in column C, then write the new device list.
# ---------------------------------------------------------------------------
$ws.Range("A29:I38").ClearContents()

$ws.Range("A29").Value = "SMA05536"
$ws.Range("B29").Value = 4.6
$ws.Range("H29").Value = "added the configuration"

$ws.Range("A30").Value = "SMA05568"
$ws.Range("B30").Value = 4.6

$ws.Range("A31").Value = "SMA05536"
$ws.Range("B31").Value = 4.6

$ws.Range("A32").Value = "SMA05545"
$ws.Range("B32").Value = 4.6

$ws.Range("A33").Value = "SMA05548"
$ws.Range("B33").Value = 4.6

$ws.Range("A34").Value = "SMA05535"
$ws.Range("B34").Value = 4.6
$ws.Range("E34").Value = "recognizing SD card as write protected - fixed"

$ws.Range("A35").Value = "SMA05619"
$ws.Range("B35").Value = 4.6

$ws.Range("A36").Value = "SMA05533"
$ws.Range("B36").Value = 4.6

$ws.Range("A37").Value = "SMA05573"
$ws.Range("B37").Value = 4.6

$ws.Range("A38").Value = "SMA05550"
$ws.Range("B38").Value = 4.6
$ws.Range("E38").Value = "tested outdoor"

# ---------------------------------------------------------------------------
# 4. Restore the on-screen selection to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("D16").Select()
